$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# Step 1: convert bsecode (column D) for rows 195-199 from text to numeric
foreach ($r in 195..199) {
    $cell = $ws.Cells.Item($r, 4)
    $numVal = [double]$cell.Value2
    $cell.Value = $numVal
}

# Step 2: append new rows 200-206

$ws.Cells.Item(200, 1).Value = 1
$ws.Cells.Item(200, 2).Value = "LTTS"
$ws.Cells.Item(200, 3).Value = "L&t Technology Services Limited"
$ws.Cells.Item(200, 4).NumberFormat = "@"
$ws.Cells.Item(200, 4).Value = "540115"
$ws.Cells.Item(200, 5).Value = 2.44
$ws.Cells.Item(200, 6).Value = 5123
$ws.Cells.Item(200, 7).Value = 267076
$ws.Cells.Item(200, 8).Value = "day"
$ws.Cells.Item(200, 9).Value = "23/07/2024 11:35:00"

$ws.Cells.Item(201, 1).Value = 2
$ws.Cells.Item(201, 2).Value = "TITAN"
$ws.Cells.Item(201, 3).Value = "Titan Company Limited"
$ws.Cells.Item(201, 4).NumberFormat = "@"
$ws.Cells.Item(201, 4).Value = "500114"
$ws.Cells.Item(201, 5).Value = 6.53
$ws.Cells.Item(201, 6).Value = 3466.85
$ws.Cells.Item(201, 7).Value = 7159507
$ws.Cells.Item(201, 8).Value = "day"
$ws.Cells.Item(201, 9).Value = "23/07/2024 11:35:00"

$ws.Cells.Item(202, 1).Value = 3
$ws.Cells.Item(202, 2).Value = "PIDILITIND"
$ws.Cells.Item(202, 3).Value = "Pidilite Industries Limited"
$ws.Cells.Item(202, 4).NumberFormat = "@"
$ws.Cells.Item(202, 4).Value = "500331"
$ws.Cells.Item(202, 5).Value = 0.83
$ws.Cells.Item(202, 6).Value = 3176.4
$ws.Cells.Item(202, 7).Value = 353111
$ws.Cells.Item(202, 8).Value = "day"
$ws.Cells.Item(202, 9).Value = "23/07/2024 11:35:00"

$ws.Cells.Item(203, 1).Value = 4
$ws.Cells.Item(203, 2).Value = "RELIANCE"
$ws.Cells.Item(203, 3).Value = "Reliance Industries Limited"
$ws.Cells.Item(203, 4).NumberFormat = "@"
$ws.Cells.Item(203, 4).Value = "500325"
$ws.Cells.Item(203, 5).Value = -0.85
$ws.Cells.Item(203, 6).Value = 2975.8
$ws.Cells.Item(203, 7).Value = 9026022
$ws.Cells.Item(203, 8).Value = "day"
$ws.Cells.Item(203, 9).Value = "23/07/2024 11:35:00"

$ws.Cells.Item(204, 1).Value = 5
$ws.Cells.Item(204, 2).Value = "INDIAMART"
$ws.Cells.Item(204, 3).Value = "Indiamart Intermesh Ltd"
$ws.Cells.Item(204, 4).NumberFormat = "@"
$ws.Cells.Item(204, 4).Value = "542726"
$ws.Cells.Item(204, 5).Value = 1.37
$ws.Cells.Item(204, 6).Value = 2895
$ws.Cells.Item(204, 7).Value = 205297
$ws.Cells.Item(204, 8).Value = "day"
$ws.Cells.Item(204, 9).Value = "23/07/2024 11:35:00"

$ws.Cells.Item(205, 1).Value = 6
$ws.Cells.Item(205, 2).Value = "BALRAMCHIN"
$ws.Cells.Item(205, 3).Value = "Balrampur Chini Mills Limited"
$ws.Cells.Item(205, 4).NumberFormat = "@"
$ws.Cells.Item(205, 4).Value = "500038"
$ws.Cells.Item(205, 5).Value = -0.11
$ws.Cells.Item(205, 6).Value = 444
$ws.Cells.Item(205, 7).Value = 1730645
$ws.Cells.Item(205, 8).Value = "day"
$ws.Cells.Item(205, 9).Value = "23/07/2024 11:35:00"

$ws.Cells.Item(206, 1).Value = 7
$ws.Cells.Item(206, 2).Value = "BIOCON"
$ws.Cells.Item(206, 3).Value = "Biocon Limited"
$ws.Cells.Item(206, 4).NumberFormat = "@"
$ws.Cells.Item(206, 4).Value = "532523"
$ws.Cells.Item(206, 5).Value = 0.51
$ws.Cells.Item(206, 6).Value = 337.05
$ws.Cells.Item(206, 7).Value = 3807462
$ws.Cells.Item(206, 8).Value = "day"
$ws.Cells.Item(206, 9).Value = "23/07/2024 11:35:00"
